# Auto-generated edit script: updates crypto price/volume table
# to match the target snapshot (commit 2023-09-06 23:07:34 UTC).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-parsed as a number by Excel's
# COM value-assignment heuristics (e.g. '215.40' -> 215.4). Force them to be
# written as literal text, matching the source inlineStr cells, then restore
# the cell style to Normal so no stray formatting diff is introduced.
function Set-TextValue($cell, $value) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = '@'
    $rng.Value = $value
    $rng.Style = 'Normal'
}

$ws.Range('D2').Value = '25.790.16'
$ws.Range('E2').Value = '  -0.30%  '
$ws.Range('D3').Value = '1.634.18'
$ws.Range('E3').Value = '  -0.16%  '
$ws.Range('E4').Value = '  -0.17%  '
Set-TextValue 'D5' '215.40'
$ws.Range('E6').Value = '  -0.84%  '
$ws.Range('E7').Value = '  -0.12%  '
$ws.Range('E8').Value = '  +0.35%  '
$ws.Range('E9').Value = '  -0.82%  '
Set-TextValue 'D10' '19.63'
$ws.Range('E10').Value = '  -3.42%  '
Set-TextValue 'D11' '0.0792'
$ws.Range('E11').Value = '  +1.47%  '
$ws.Range('E12').Value = '  -0.13%  '
$ws.Range('B13').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C13').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D13').Value = '1.859.84'
$ws.Range('E13').Value = '  -0.15%  '
$ws.Range('B14').Value = 'WrappedEther'
$ws.Range('C14').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D14').Value = '1.634.51'
$ws.Range('E14').Value = '  -0.14%  '
$ws.Range('E15').Value = '  -0.18%  '
$ws.Range('D16').Value = '0.0₃0767'
$ws.Range('E16').Value = '  -0.06%  '
Set-TextValue 'D17' '62.79'
$ws.Range('E17').Value = '  -0.71%  '
$ws.Range('D18').Value = '25.796.60'
$ws.Range('E18').Value = '  -0.30%  '
$ws.Range('E19').Value = '  -0.17%  '
$ws.Range('E20').Value = '  +1.54%  '
Set-TextValue 'D21' '194.37'
$ws.Range('E21').Value = '  +0.14%  '
$ws.Range('E22').Value = '  +0.27%  '
$ws.Range('E23').Value = '  +1.69%  '
$ws.Range('E24').Value = '  -0.09%  '
$ws.Range('E25').Value = '  +2.12%  '
Set-TextValue 'D26' '142.74'
$ws.Range('E26').Value = '  +2.92%  '
$ws.Range('E27').Value = '  +0.01%  '
Set-TextValue 'D29' '15.57'
$ws.Range('E29').Value = '  -0.12%  '
$ws.Range('E30').Value = '  -0.24%  '
Set-TextValue 'D31' '0.0494'
$ws.Range('E31').Value = '  -0.45%  '
$ws.Range('E32').Value = '  +1.50%  '
$ws.Range('E33').Value = '  -0.26%  '
$ws.Range('E34').Value = '  +0.78%  '
Set-TextValue 'D35' '2.38'
$ws.Range('E35').Value = '  -0.06%  '
Set-TextValue 'D36' '0.905'
$ws.Range('E36').Value = '  +0.10%  '
$ws.Range('D37').Value = '1.131.06'
$ws.Range('E37').Value = '  -0.39%  '
Set-TextValue 'D39' '0.547'
$ws.Range('E39').Value = '  -1.77%  '
Set-TextValue 'D40' '0.0156'
$ws.Range('E40').Value = '  -0.35%  '
$ws.Range('E41').Value = '  -0.26%  '
Set-TextValue 'D43' '100.49'
$ws.Range('E43').Value = '  +1.19%  '
$ws.Range('E44').Value = '  +0.80%  '
$ws.Range('D45').Value = '1.768.76'
$ws.Range('E45').Value = '  -0.37%  '
$ws.Range('B46').Value = 'BabyDogeCoin'
$ws.Range('C46').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D46').Value = '0.0₆0112'
$ws.Range('E46').Value = '  +0.97%  '
$ws.Range('B47').Value = 'Aave'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue 'D47' '55.24'
$ws.Range('E47').Value = '  -0.46%  '
$ws.Range('B48').Value = 'Cronos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue 'D48' '0.0506'
$ws.Range('E48').Value = '  -0.31%  '
$ws.Range('B49').Value = 'Mantle'
$ws.Range('C49').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextValue 'D49' '0.417'
$ws.Range('E49').Value = '  -2.29%  '
$ws.Range('E50').Value = '  +0.36%  '
$ws.Range('B51').Value = 'EnergySwap'
$ws.Range('C51').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue 'D51' '7.52'
$ws.Range('E51').Value = '  -3.56%  '
